$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 330.8
$ws.Cells.Item(4, 9).Value = 57
$ws.Cells.Item(4, 10).Value = 513.3333
$ws.Cells.Item(4, 11).Value = 57
$ws.Cells.Item(4, 12).Value = 513.3333
$ws.Cells.Item(4, 13).Value = 57
$ws.Cells.Item(4, 14).Value = -741.3333
$ws.Cells.Item(17, 8).Value = 1559438
$ws.Cells.Item(17, 10).Value = 1590593.5
$ws.Cells.Item(17, 12).Value = 4771780.5
$ws.Cells.Item(17, 14).Value = -4772116.5
$ws.Cells.Item(28, 8).Value = 913.5454999999999
$ws.Cells.Item(28, 10).Value = 2501.3333
$ws.Cells.Item(28, 12).Value = 2501.3333
$ws.Cells.Item(28, 14).Value = -3471.3333
$ws.Cells.Item(62, 8).Value = 6791.7334
$ws.Cells.Item(62, 9).Value = 6603.909
$ws.Cells.Item(62, 10).Value = 7308.25
$ws.Cells.Item(62, 11).Value = 6603.909
$ws.Cells.Item(62, 12).Value = 7308.25
$ws.Cells.Item(62, 13).Value = -5979.909
$ws.Cells.Item(62, 14).Value = -8556.25
$ws.Cells.Item(65, 8).Value = 6791.7334
$ws.Cells.Item(65, 9).Value = 6603.909
$ws.Cells.Item(65, 10).Value = 7308.25
$ws.Cells.Item(65, 11).Value = 33019.545
$ws.Cells.Item(65, 12).Value = 36541.25
$ws.Cells.Item(65, 13).Value = -29899.545
$ws.Cells.Item(65, 14).Value = -42781.25
$ws.Cells.Item(74, 8).Value = 7458.1665
$ws.Cells.Item(74, 10).Value = 9099.6
$ws.Cells.Item(74, 12).Value = 9099.6
$ws.Cells.Item(74, 14).Value = -10971.6
$ws.Cells.Item(77, 8).Value = 7458.1665
$ws.Cells.Item(77, 10).Value = 9099.6
$ws.Cells.Item(77, 12).Value = 45498
$ws.Cells.Item(77, 14).Value = -54858
$ws.Cells.Item(103, 8).Value = 369.5
$ws.Cells.Item(103, 9).Value = 260
$ws.Cells.Item(103, 11).Value = 780
$ws.Cells.Item(103, 13).Value = -194
$ws.Cells.Item(113, 8).Value = 3331.125
$ws.Cells.Item(113, 10).Value = 3841.5
$ws.Cells.Item(113, 12).Value = 3841.5
$ws.Cells.Item(113, 14).Value = -10349.5
$ws.Cells.Item(116, 8).Value = 3338.9473
$ws.Cells.Item(116, 10).Value = 2882.6667
$ws.Cells.Item(116, 12).Value = 2882.6667
$ws.Cells.Item(116, 14).Value = -9766.6667
$ws.Cells.Item(132, 8).Value = 9074.868
$ws.Cells.Item(132, 9).Value = 3083.5293
$ws.Cells.Item(132, 11).Value = 9250.5879
$ws.Cells.Item(132, 13).Value = -6720.5879

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 7842.6816
$ws.Cells.Item(32, 9).Value = 7842.6816
$ws.Cells.Item(32, 11).Value = 7842.6816
$ws.Cells.Item(32, 13).Value = -7555.6816
$ws.Cells.Item(45, 8).Value = 4106.778
$ws.Cells.Item(45, 9).Value = 3740.5
$ws.Cells.Item(45, 10).Value = 4399.8
$ws.Cells.Item(45, 11).Value = 3740.5
$ws.Cells.Item(45, 12).Value = 4399.8
$ws.Cells.Item(45, 13).Value = -3363.5
$ws.Cells.Item(45, 14).Value = -5153.8
$ws.Cells.Item(97, 8).Value = 1386.3334
$ws.Cells.Item(97, 9).Value = 1386.3334
$ws.Cells.Item(97, 11).Value = 1386.3334
$ws.Cells.Item(97, 13).Value = -890.3334
$ws.Cells.Item(102, 8).Value = 2801.3684
$ws.Cells.Item(102, 9).Value = 1772.1
$ws.Cells.Item(102, 10).Value = 3945
$ws.Cells.Item(102, 11).Value = 1772.1
$ws.Cells.Item(102, 12).Value = 3945
$ws.Cells.Item(102, 13).Value = -150.0999999999999
$ws.Cells.Item(102, 14).Value = -7189
$ws.Cells.Item(110, 8).Value = 3654.3333
$ws.Cells.Item(110, 10).Value = 4003.75
$ws.Cells.Item(110, 12).Value = 4003.75
$ws.Cells.Item(110, 14).Value = -8093.75
$ws.Cells.Item(122, 8).Value = 2339.3103
$ws.Cells.Item(122, 9).Value = 1994.1482
$ws.Cells.Item(122, 10).Value = 6999
$ws.Cells.Item(122, 11).Value = 5982.444600000001
$ws.Cells.Item(122, 12).Value = 20997
$ws.Cells.Item(122, 13).Value = -3532.444600000001
$ws.Cells.Item(122, 14).Value = -25897

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 4842.25
$ws.Cells.Item(99, 10).Value = 7999.6665
$ws.Cells.Item(99, 12).Value = 7999.6665
$ws.Cells.Item(99, 14).Value = -10995.6665
$ws.Cells.Item(107, 8).Value = 35715544
$ws.Cells.Item(107, 9).Value = 50001096
$ws.Cells.Item(107, 11).Value = 50001096
$ws.Cells.Item(107, 13).Value = -49999176

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3084.5667
$ws.Cells.Item(31, 9).Value = 1474.3334
$ws.Cells.Item(31, 10).Value = 5499.9165
$ws.Cells.Item(31, 11).Value = 1474.3334
$ws.Cells.Item(31, 12).Value = 5499.9165
$ws.Cells.Item(31, 13).Value = -1179.3334
$ws.Cells.Item(31, 14).Value = -6089.9165
$ws.Cells.Item(34, 8).Value = 3084.5667
$ws.Cells.Item(34, 9).Value = 1474.3334
$ws.Cells.Item(34, 10).Value = 5499.9165
$ws.Cells.Item(34, 11).Value = 1474.3334
$ws.Cells.Item(34, 12).Value = 5499.9165
$ws.Cells.Item(34, 13).Value = -1272.3334
$ws.Cells.Item(34, 14).Value = -5903.9165
$ws.Cells.Item(99, 8).Value = 2770.353
$ws.Cells.Item(99, 10).Value = 3499
$ws.Cells.Item(99, 12).Value = 3499
$ws.Cells.Item(99, 14).Value = -6495
$ws.Cells.Item(105, 8).Value = 2019.8462
$ws.Cells.Item(105, 9).Value = 2059.8333
$ws.Cells.Item(105, 11).Value = 2059.8333
$ws.Cells.Item(105, 13).Value = -312.8332999999998
$ws.Cells.Item(106, 8).Value = 78221.336
$ws.Cells.Item(106, 10).Value = 78221.336
$ws.Cells.Item(106, 12).Value = 78221.336
$ws.Cells.Item(106, 14).Value = -80745.336
$ws.Cells.Item(107, 8).Value = 1404.5
$ws.Cells.Item(107, 9).Value = 557.4286
$ws.Cells.Item(107, 11).Value = 557.4286
$ws.Cells.Item(107, 13).Value = 1362.5714
$ws.Cells.Item(126, 8).Value = 2770.353
$ws.Cells.Item(126, 10).Value = 3499
$ws.Cells.Item(126, 12).Value = 10497
$ws.Cells.Item(126, 14).Value = -15437
$ws.Cells.Item(132, 8).Value = 1820793.1
$ws.Cells.Item(132, 9).Value = 2107866
$ws.Cells.Item(132, 11).Value = 6323598
$ws.Cells.Item(132, 13).Value = -6321068
$ws.Cells.Item(134, 8).Value = 2785.725
$ws.Cells.Item(134, 9).Value = 1385.0834
$ws.Cells.Item(134, 11).Value = 4155.2502
$ws.Cells.Item(134, 13).Value = -1620.2502

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 1522.4667
$ws.Cells.Item(113, 10).Value = 1634
$ws.Cells.Item(113, 12).Value = 4902
$ws.Cells.Item(113, 14).Value = -9242
$ws.Cells.Item(121, 8).Value = 3740.2727
$ws.Cells.Item(121, 9).Value = 447.25
$ws.Cells.Item(121, 10).Value = 5622
$ws.Cells.Item(121, 11).Value = 1341.75
$ws.Cells.Item(121, 12).Value = 16866
$ws.Cells.Item(121, 13).Value = -31.75
$ws.Cells.Item(121, 14).Value = -19486

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 2229.0557
$ws.Cells.Item(97, 9).Value = 1510.5
$ws.Cells.Item(97, 10).Value = 3666.1667
$ws.Cells.Item(97, 11).Value = 1510.5
$ws.Cells.Item(97, 12).Value = 3666.1667
$ws.Cells.Item(97, 13).Value = -1014.5
$ws.Cells.Item(97, 14).Value = -4658.1667
$ws.Cells.Item(102, 8).Value = 26652.625
$ws.Cells.Item(102, 9).Value = 1338.4839
$ws.Cells.Item(102, 11).Value = 1338.4839
$ws.Cells.Item(102, 13).Value = 283.5161000000001
$ws.Cells.Item(107, 8).Value = 817.5909
$ws.Cells.Item(107, 9).Value = 585.2727
$ws.Cells.Item(107, 10).Value = 1049.909
$ws.Cells.Item(107, 11).Value = 585.2727
$ws.Cells.Item(107, 12).Value = 1049.909
$ws.Cells.Item(107, 13).Value = 1334.7273
$ws.Cells.Item(107, 14).Value = -4889.909
$ws.Cells.Item(126, 8).Value = 3896.5715
$ws.Cells.Item(126, 9).Value = 4208
$ws.Cells.Item(126, 10).Value = 3481.3333
$ws.Cells.Item(126, 11).Value = 12624
$ws.Cells.Item(126, 12).Value = 10443.9999
$ws.Cells.Item(126, 13).Value = -10154
$ws.Cells.Item(126, 14).Value = -15383.9999
$ws.Cells.Item(132, 8).Value = 14496252
$ws.Cells.Item(132, 9).Value = 16669917
$ws.Cells.Item(132, 10).Value = 5154.6665
$ws.Cells.Item(132, 11).Value = 50009751
$ws.Cells.Item(132, 12).Value = 15463.9995
$ws.Cells.Item(132, 13).Value = -50007221
$ws.Cells.Item(132, 14).Value = -20523.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(42, 8).Value = 27156.076
$ws.Cells.Item(42, 10).Value = 22201
$ws.Cells.Item(42, 12).Value = 22201
$ws.Cells.Item(42, 14).Value = -23327
$ws.Cells.Item(49, 8).Value = 27156.076
$ws.Cells.Item(49, 10).Value = 22201
$ws.Cells.Item(49, 12).Value = 22201
$ws.Cells.Item(49, 14).Value = -22495
$ws.Cells.Item(93, 8).Value = 1449.4546
$ws.Cells.Item(93, 9).Value = 1660.6666
$ws.Cells.Item(93, 11).Value = 1660.6666
$ws.Cells.Item(93, 13).Value = -412.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 589.1070999999999
$ws.Cells.Item(107, 9).Value = 289.41177
$ws.Cells.Item(107, 11).Value = 868.23531
$ws.Cells.Item(107, 13).Value = 1051.76469
